# Ajout des colonnes aux matrices
# Adds a header row (Hostname / Conf réseau IP / Technologies réseaux actives)
# to the first sheet, styles it (bold, Arial, 10pt) and sets up column widths
# + a few formatting-only cells to the right, matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header values -------------------------------------------------------
$ws.Range("A1").Value = "Hostname"
$ws.Range("B1").Value = "Conf réseau IP"
$ws.Range("C1").Value = "Technologies réseaux actives"

# --- Header style: bold, Arial, 10pt -------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Name = "Arial"

# Propagate the exact same style to the rest of the header row (B1:E1 and
# G1:H1 -- F1 is intentionally left untouched) without creating additional
# style/font definitions.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:E1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 27.333333333333332

# --- Selection --------------------------------------------------------------
$ws.Range("C10").Select() | Out-Null
